$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 7441.4287
$ws.Range("I62").Value = 7441.4287
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 7441.4287
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -6817.4287
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 7441.4287
$ws.Range("I65").Value = 7441.4287
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 37207.14350000001
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -34087.14350000001
$ws.Range("N65").ClearContents()
# Row 86
$ws.Range("H86").Value = 1383.2142
$ws.Range("I86").Value = 1280.6666
$ws.Range("K86").Value = 1280.6666
$ws.Range("M86").Value = -157.6666
# Row 89
$ws.Range("H89").Value = 1383.2142
$ws.Range("I89").Value = 1280.6666
$ws.Range("K89").Value = 6403.333000000001
$ws.Range("M89").Value = -787.3330000000005
# Row 138
$ws.Range("H138").Value = 1625.84
$ws.Range("J138").Value = 2054.4915
$ws.Range("L138").Value = 6163.4745
$ws.Range("N138").Value = -16443.4745

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1111868.6
$ws.Range("I2").Value = 1111868.6
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1111868.6
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1111755.6
$ws.Range("N2").ClearContents()
# Row 45
$ws.Range("H45").Value = 1710.2222
$ws.Range("I45").Value = 1399.5
$ws.Range("K45").Value = 1399.5
$ws.Range("M45").Value = -1022.5
# Row 63
$ws.Range("H63").Value = 2184.1667
$ws.Range("I63").Value = 1727.5
$ws.Range("J63").Value = 2412.5
$ws.Range("K63").Value = 1727.5
$ws.Range("L63").Value = 2412.5
$ws.Range("M63").Value = -1041.5
$ws.Range("N63").Value = -3784.5
# Row 66
$ws.Range("H66").Value = 2184.1667
$ws.Range("I66").Value = 1727.5
$ws.Range("J66").Value = 2412.5
$ws.Range("K66").Value = 8637.5
$ws.Range("L66").Value = 12062.5
$ws.Range("M66").Value = -5205.5
$ws.Range("N66").Value = -18926.5
# Row 116
$ws.Range("H116").Value = 1111868.6
$ws.Range("I116").Value = 1111868.6
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1111868.6
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1109574.6
$ws.Range("N116").ClearContents()
# Row 123
$ws.Range("H123").Value = 81992.5
$ws.Range("J123").Value = 81992.5
$ws.Range("L123").Value = 81992.5
$ws.Range("N123").Value = -91792.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1111868.6
$ws.Range("I3").Value = 1111868.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1111868.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1111754.6
$ws.Range("N3").ClearContents()
# Row 86
$ws.Range("H86").Value = 1820.2
$ws.Range("J86").Value = 1703.6666
$ws.Range("L86").Value = 1703.6666
$ws.Range("N86").Value = -3949.6666
# Row 89
$ws.Range("H89").Value = 1820.2
$ws.Range("J89").Value = 1703.6666
$ws.Range("L89").Value = 8518.333000000001
$ws.Range("N89").Value = -19750.333
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 107
$ws.Range("H107").Value = 740.5
$ws.Range("I107").Value = 649.6667
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 649.6667
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = 1270.3333
$ws.Range("N107").Value = -4853
# Row 134
$ws.Range("H134").Value = 5841.5586
$ws.Range("I134").Value = 7358.1304
$ws.Range("K134").Value = 22074.3912
$ws.Range("M134").Value = -19539.3912

$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 420
# Row 62
$ws.Range("H62").Value = 3429
$ws.Range("I62").Value = 2904
$ws.Range("J62").Value = 3604
$ws.Range("K62").Value = 2904
$ws.Range("L62").Value = 3604
$ws.Range("M62").Value = -2280
$ws.Range("N62").Value = -4852
# Row 65
$ws.Range("H65").Value = 3429
$ws.Range("I65").Value = 2904
$ws.Range("J65").Value = 3604
$ws.Range("K65").Value = 14520
$ws.Range("L65").Value = 18020
$ws.Range("M65").Value = -11400
$ws.Range("N65").Value = -24260
# Row 107
$ws.Range("H107").Value = 588.4400000000001
$ws.Range("J107").Value = 800
$ws.Range("L107").Value = 800
$ws.Range("N107").Value = -4640

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 524.375
# Row 50
$ws.Range("H50").Value = 100070984
$ws.Range("J50").Value = 166668350
$ws.Range("L50").Value = 500005050
$ws.Range("N50").Value = -500006012
# Row 53
$ws.Range("H53").Value = 100070984
$ws.Range("J53").Value = 166668350
$ws.Range("L53").Value = 500005050
$ws.Range("N53").Value = -500006012
# Row 122
$ws.Range("H122").Value = 1006.38464
$ws.Range("I122").Value = 749.5
$ws.Range("K122").Value = 6745.5
$ws.Range("M122").Value = -4295.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4445.0557
$ws.Range("I70").Value = 4206.727
$ws.Range("J70").Value = 4819.5713
$ws.Range("K70").Value = 4206.727
$ws.Range("L70").Value = 4819.5713
$ws.Range("M70").Value = -3936.727
$ws.Range("N70").Value = -5359.5713
# Row 73
$ws.Range("H73").Value = 4445.0557
$ws.Range("I73").Value = 4206.727
$ws.Range("J73").Value = 4819.5713
$ws.Range("K73").Value = 4206.727
$ws.Range("L73").Value = 4819.5713
$ws.Range("M73").Value = -3270.727
$ws.Range("N73").Value = -6691.5713
# Row 107
$ws.Range("H107").Value = 250
$ws.Range("I107").Value = 250
$ws.Range("K107").Value = 250
$ws.Range("M107").Value = 1670
# Row 113
$ws.Range("H113").Value = 1995
$ws.Range("J113").Value = 1995
$ws.Range("L113").Value = 1995
$ws.Range("N113").Value = -6335
# Row 122
$ws.Range("H122").Value = 1491.1538
$ws.Range("I122").Value = 1441.2858
$ws.Range("J122").Value = 1549.3334
$ws.Range("K122").Value = 4323.857400000001
$ws.Range("L122").Value = 4648.0002
$ws.Range("M122").Value = -1873.857400000001
$ws.Range("N122").Value = -9548.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()
# Row 68
$ws.Range("H68").Value = 2411
$ws.Range("I68").Value = 1660
$ws.Range("K68").Value = 1660
$ws.Range("M68").Value = -911
# Row 71
$ws.Range("H71").Value = 2411
$ws.Range("I71").Value = 1660
$ws.Range("K71").Value = 8300
$ws.Range("M71").Value = -4556
# Row 104
$ws.Range("H104").Value = 107500
$ws.Range("J104").Value = 107500
$ws.Range("L104").Value = 107500
$ws.Range("N104").Value = -114488
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 132
$ws.Range("H132").Value = 2046.9302
$ws.Range("I132").Value = 1744.1428
$ws.Range("K132").Value = 5232.428400000001
$ws.Range("M132").Value = -2702.428400000001
# Row 136
$ws.Range("H136").Value = 3135.3333
$ws.Range("I136").Value = 3173.2856
$ws.Range("J136").Value = 3002.5
$ws.Range("K136").Value = 9519.856800000001
$ws.Range("L136").Value = 9007.5
$ws.Range("M136").Value = -6969.856800000001
$ws.Range("N136").Value = -14107.5

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 28670.4
$ws.Range("J16").Value = 28670.4
$ws.Range("L16").Value = 28670.4
$ws.Range("N16").Value = -29254.4
# Row 104
$ws.Range("H104").Value = 29999.5
$ws.Range("J104").Value = 29999.5
$ws.Range("L104").Value = 29999.5
$ws.Range("N104").Value = -36987.5
